$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.682.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.88%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.920.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.43%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'335.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.48%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.05%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = "'0.4111"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.29%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'48.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.08%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.08038"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.37%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.018"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.06%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +5.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.928.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.02%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.000"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.17%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +3.67%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'90.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.02%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.9996"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.16%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +2.36%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06592"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.18%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +5.21%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.00%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'29.634.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.73%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.584"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +6.16%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +10.96%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.209"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.75%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.171.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +4.87%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'156.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.08%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +4.03%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +5.69%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.734"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +9.27%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'117.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'1.069"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +14.74%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09483"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.67%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.438"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.79%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.574"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.67%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.425"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +4.95%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'Hedera"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.06145"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.60%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'VeChain"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.02274"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.29%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'8.446"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.68%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.183"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.96%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.5894"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.64%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1848"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.89%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'10.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.15%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.261"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.14%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.358"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.78%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.07507"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.54%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Decentraland"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.5589"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.55%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'12.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.90%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +4.22%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'113.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.83%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.3000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +14.47%  "
$ws.Range("E51").Style = "Normal"
